# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colour scheme (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" colour scheme (used by the Slide Master
#                             and therefore by every slide in the presentation)
#
# The authored edit swaps the content of the two theme parts, which - from the
# presentation's point of view - means the design actually applied to the slides
# switches from the "Integral" / "Red Violet" palette to the plain "Office" palette.
# Re-colour the active theme (the one driving the Slide Master / all slides) to the
# Office theme's standard 12-colour palette so the deck's design matches that swap.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Theme colour order exposed by ThemeColorScheme.Item(n):
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5-10 Accent1-Accent6, 11 Hyperlink, 12 Followed Hyperlink
# Target values are the Office theme's colours (RRGGBB), converted to the
# BGR-packed OLE_COLOR integers (R | G<<8 | B<<16) expected by .RGB.

$colorScheme.Item(1).RGB  = 0          # Dark1          000000
$colorScheme.Item(2).RGB  = 16777215   # Light1         FFFFFF
$colorScheme.Item(3).RGB  = 6968388    # Dark2          44546A
$colorScheme.Item(4).RGB  = 15132391   # Light2         E7E6E6
$colorScheme.Item(5).RGB  = 13998939   # Accent1        5B9BD5
$colorScheme.Item(6).RGB  = 3243501    # Accent2        ED7D31
$colorScheme.Item(7).RGB  = 10855845   # Accent3        A5A5A5
$colorScheme.Item(8).RGB  = 49407      # Accent4        FFC000
$colorScheme.Item(9).RGB  = 12874308   # Accent5        4472C4
$colorScheme.Item(10).RGB = 4697456    # Accent6        70AD47
$colorScheme.Item(11).RGB = 12673797   # Hyperlink      0563C1
$colorScheme.Item(12).RGB = 7491477    # FollowedHyperlink 954F72
